# The commit adds three new daily price records (Murcott / Mandarina) that
# belong chronologically right before the existing "Clemenuless 44391" block
# that starts at row 185. Excel users would do this by inserting three blank
# rows above row 185 (shifting every following row down by three, 273 -> 276
# total data rows) and then typing in the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 185, pushing old rows 185:273 down
# to 188:276.
$ws.Rows.Item(185).EntireRow.Insert()
$ws.Rows.Item(185).EntireRow.Insert()
$ws.Rows.Item(185).EntireRow.Insert()

# Common/fixed values for every data row in this sheet.
$mercadoId   = 2
$mercado     = "Comercializadora del Agro de Limarí"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100102
$producto    = "Cítricos"
$categoriaId = 100102004
$categoria   = "Mandarina"
$origen      = "Provincia de Limarí"
$unidad      = "$/bandeja 10 kilos"

# Row 185: Murcott / Especial
$r = 185
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 44489
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Murcott"
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 500
$ws.Cells.Item($r, 14).Value = 4500
$ws.Cells.Item($r, 15).Value = 5000
$ws.Cells.Item($r, 16).Value = 4750
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 475
$ws.Cells.Item($r, 20).Value = 10

# Row 186: Murcott / Primera
$r = 186
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 44489
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Murcott"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 500
$ws.Cells.Item($r, 14).Value = 3500
$ws.Cells.Item($r, 15).Value = 4000
$ws.Cells.Item($r, 16).Value = 3750
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 375
$ws.Cells.Item($r, 20).Value = 10

# Row 187: Murcott / Segunda
$r = 187
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 44489
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Murcott"
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 400
$ws.Cells.Item($r, 14).Value = 2500
$ws.Cells.Item($r, 15).Value = 3000
$ws.Cells.Item($r, 16).Value = 2750
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 275
$ws.Cells.Item($r, 20).Value = 10
